$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($addr, $val)
    $escaped = $val -replace '"', '""'
    $ws.Range($addr).Formula = '="' + $escaped + '"'
    $ws.Range($addr).Copy($null)
    $ws.Range($addr).PasteSpecial(-4163)
}

Set-TextValue "D2" '62.713.03'
Set-TextValue "E2" '  -2.00%  '

Set-TextValue "D3" '2.579.37'
Set-TextValue "E3" '  -1.96%  '

Set-TextValue "E4" '  +0.00%  '

Set-TextValue "D5" '567.98'
Set-TextValue "E5" '  -1.63%  '

Set-TextValue "D6" '153.11'
Set-TextValue "E6" '  -2.43%  '

Set-TextValue "E7" '  +0.02%  '

Set-TextValue "D8" '0.615'
Set-TextValue "E8" '  -2.56%  '

Set-TextValue "E9" '  -4.94%  '

Set-TextValue "D10" '5.69'
Set-TextValue "E10" '  -2.22%  '

Set-TextValue "B11" 'TRON'
Set-TextValue "C11" 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
Set-TextValue "D11" '0.156'
Set-TextValue "E11" '  +0.24%  '

Set-TextValue "B12" 'Cardano'
Set-TextValue "C12" 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
Set-TextValue "D12" '0.374'
Set-TextValue "E12" '  -2.90%  '

Set-TextValue "D13" '27.87'
Set-TextValue "E13" '  -2.16%  '

Set-TextValue "D14" '3.054.09'
Set-TextValue "E14" '  -1.74%  '

Set-TextValue "E15" '  -2.88%  '

Set-TextValue "D16" '62.649.44'
Set-TextValue "E16" '  -1.85%  '

Set-TextValue "D17" '2.606.67'
Set-TextValue "E17" '  -1.02%  '

Set-TextValue "D18" '11.84'
Set-TextValue "E18" '  -2.48%  '

Set-TextValue "D19" '7.44'
Set-TextValue "E19" '  -3.38%  '

Set-TextValue "D20" '4.40'
Set-TextValue "E20" '  -4.49%  '

Set-TextValue "D21" '335.20'
Set-TextValue "E21" '  -2.87%  '

Set-TextValue "E22" '  +0.07%  '

Set-TextValue "D23" '67.07'
Set-TextValue "E23" '  -0.50%  '

Set-TextValue "D24" '1.84'
Set-TextValue "E24" '  +5.00%  '

Set-TextValue "E25" '  -1.12%  '

Set-TextValue "D26" '1.60'
Set-TextValue "E26" '  +1.07%  '

Set-TextValue "D27" '8.99'
Set-TextValue "E27" '  -2.88%  '

Set-TextValue "D28" '560.44'
Set-TextValue "E28" '  -6.42%  '

Set-TextValue "D29" '8.00'
Set-TextValue "E29" '  +1.11%  '

Set-TextValue "E30" '  +0.38%  '

Set-TextValue "D31" '0.157'
Set-TextValue "E31" '  -3.08%  '

Set-TextValue "E32" '  -4.33%  '

Set-TextValue "E33" '  -2.84%  '

Set-TextValue "D34" '6.39'
Set-TextValue "E34" '  -3.56%  '

Set-TextValue "D35" '5.25'
Set-TextValue "E35" '  -1.58%  '

Set-TextValue "E36" '  +0.29%  '

Set-TextValue "D37" '0.394'
Set-TextValue "E37" '  -3.64%  '

Set-TextValue "D38" '19.31'
Set-TextValue "E38" '  -3.09%  '

Set-TextValue "D39" '153.78'
Set-TextValue "E39" '  -0.71%  '

Set-TextValue "E40" '  -2.42%  '

Set-TextValue "D42" '2.48'
Set-TextValue "E42" '  +2.29%  '

Set-TextValue "D43" '157.68'
Set-TextValue "E43" '  +0.27%  '

Set-TextValue "D44" '23.50'
Set-TextValue "E44" '  +1.20%  '

Set-TextValue "D45" '3.82'
Set-TextValue "E45" '  -3.33%  '

Set-TextValue "E46" '  -3.99%  '

Set-TextValue "E47" '  -1.70%  '

Set-TextValue "D48" '0.0983'
Set-TextValue "E48" '  -3.78%  '

Set-TextValue "D49" '0.0240'
Set-TextValue "E49" '  -4.19%  '

Set-TextValue "D50" '0.0₆0229'
Set-TextValue "E50" '  -2.25%  '

Set-TextValue "E51" '  -1.35%  '

$excel.CutCopyMode = $false